$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 2919.202174992006
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 2459690191846.092
$ws.Range("G2").Value = 2459690194768.559
